$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.386.86"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.680.10"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'315.96"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'0.3879"
$ws.Range("E7").Value = "  -1.36%  "
$ws.Range("D8").Value = "'0.3994"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").Value = "'1.475"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").Value = "'52.19"
$ws.Range("E11").Value = "  -3.42%  "
$ws.Range("D12").Value = "'0.08726"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  +10.05%  "
$ws.Range("D14").Value = "'7.456"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").Value = "'7.963"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "'0.00001340"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").Value = "1.671.50"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "'97.67"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").Value = "'0.07205"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").Value = "'19.58"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "'7.236"
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'14.13"
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("D24").Value = "24.359.74"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "'3.003"
$ws.Range("E25").Value = "  -6.63%  "
$ws.Range("D26").Value = "'2.337"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("D28").Value = "'167.07"
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("D29").Value = "'8.615"
$ws.Range("E29").Value = "  +11.06%  "
$ws.Range("D30").Value = "'5.346"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("D31").Value = "'137.99"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").Value = "1.855.16"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "'0.08756"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'7.348"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("D35").Value = "'1.044"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").Value = "'1.968"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "'0.02942"
$ws.Range("E37").Value = "  +7.57%  "
$ws.Range("D38").Value = "'0.2743"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").Value = "'10.75"
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("D40").Value = "'0.09131"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").Value = "'13.99"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").Value = "'0.7907"
$ws.Range("E42").Value = "  +2.76%  "
$ws.Range("D43").Value = "'1.471"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "'17.28"
$ws.Range("E44").Value = "  +7.91%  "
$ws.Range("D45").Value = "'0.7179"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'2.576"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "'4.257"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "'1.404"
$ws.Range("E48").Value = "  +7.28%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'139.24"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").Value = "'0.08026"
$ws.Range("E51").Value = "  +0.44%  "
